# The workbook tracks weekly price observations for "Acelga" at the
# "Femacal de La Calera" market. This edit adds one new weekly
# observation, inserted as a new row right before the existing row 125
# (pushing the existing row 125 and everything below it down by one row,
# so the former last row, 189, becomes row 190).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 125; existing rows 125-189 shift to 126-190.
$ws.Rows.Item(125).EntireRow.Insert()

# Populate the newly inserted row 125 with the new observation.
$ws.Range("A125").Value = 3
$ws.Range("B125").Value = "Femacal de La Calera"
$ws.Range("C125").Value = "Coquimbo"
$ws.Range("D125").Value = 44452
$ws.Range("E125").Value = 5
$ws.Range("F125").Value = 100112009
$ws.Range("G125").Value = "Acelga"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 270
$ws.Range("K125").Value = 2000
$ws.Range("L125").Value = 2200
$ws.Range("M125").Value = 2089
$ws.Range("N125").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O125").Value = "Provincia de Quillota"
$ws.Range("P125").Value = 348
$ws.Range("Q125").Value = 6
$ws.Range("R125").Value = "Hortaliza"
